$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")  # "config" is the first/active sheet

# New strings are appended to the shared string table in the order the
# values are first assigned, so the order below matches the authoring
# order: dBm, Psat input, Pout, Overdrive input, ATR start input dBm.
$ws.Range("L2").Value = "dBm"
$ws.Range("M1").Value = "Psat input"
$ws.Range("O1").Value = "Pout"
$ws.Range("N1").Value = "Overdrive input"
$ws.Range("L1").Value = "ATR start input dBm"

# Remaining unit row (row 2) values for columns M:O reuse the "dBm" string
$ws.Range("M2").Value = "dBm"
$ws.Range("N2").Value = "dBm"
$ws.Range("O2").Value = "dBm"

# New data row (row 3) values for columns L:O
$ws.Range("L3").Value = -20
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 50

# Updated data row (row 3) existing columns
$ws.Range("H3").Value = 18
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 16

# Column widths for the new columns (target OOXML <col width="..."/> is
# 18.59765625 / 9.3984375 / 14.3984375 - the engine's ColumnWidth setter
# snaps the stored width to the nearest 1/7 character unit, so these are
# the closest values it can represent).
$ws.Columns.Item(12).ColumnWidth = 17.857142857142858
$ws.Columns.Item(13).ColumnWidth = 8.714285714285714
$ws.Columns.Item(14).ColumnWidth = 13.714285714285714

# Update selection to reflect new active cell
$ws.Range("H6").Select()
